$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 799.3333
$ws.Range("I4").Value = 799.3333
$ws.Range("K4").Value = 799.3333
$ws.Range("M4").Value = -685.3333
$ws.Range("H18").Value = 250003000
$ws.Range("I18").Value = 3000
$ws.Range("K18").Value = 3000
$ws.Range("M18").Value = -2716
$ws.Range("H29").Value = 6479.5
$ws.Range("J29").Value = 6479.5
$ws.Range("L29").Value = 19438.5
$ws.Range("N29").Value = -20000.5
$ws.Range("H98").Value = 2433.6667
$ws.Range("I98").Value = 2594.76
$ws.Range("J98").Value = 420
$ws.Range("K98").Value = 2594.76
$ws.Range("L98").Value = 420
$ws.Range("M98").Value = -1096.76
$ws.Range("N98").Value = -3416
$ws.Range("H122").Value = 2433.6667
$ws.Range("I122").Value = 2594.76
$ws.Range("J122").Value = 420
$ws.Range("K122").Value = 7784.280000000001
$ws.Range("L122").Value = 1260
$ws.Range("M122").Value = -5334.280000000001
$ws.Range("N122").Value = -6160
$ws.Range("H138").Value = 3004.0186
$ws.Range("J138").Value = 2598.5
$ws.Range("L138").Value = 7795.5
$ws.Range("N138").Value = -18075.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2867.7188
$ws.Range("I32").Value = 2465.8772
$ws.Range("K32").Value = 2465.8772
$ws.Range("M32").Value = -2178.8772
$ws.Range("H45").Value = 40090.363
$ws.Range("I45").Value = 43498
$ws.Range("K45").Value = 43498
$ws.Range("M45").Value = -43121
$ws.Range("H74").Value = 294465.25
$ws.Range("I74").Value = 506147.62
$ws.Range("J74").Value = 3402
$ws.Range("K74").Value = 506147.62
$ws.Range("L74").Value = 3402
$ws.Range("M74").Value = -505273.62
$ws.Range("N74").Value = -5150
$ws.Range("H77").Value = 294465.25
$ws.Range("I77").Value = 506147.62
$ws.Range("J77").Value = 3402
$ws.Range("K77").Value = 2530738.1
$ws.Range("L77").Value = 17010
$ws.Range("M77").Value = -2526370.1
$ws.Range("N77").Value = -25746
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 1455.35
$ws.Range("I122").Value = 1455.35
$ws.Range("K122").Value = 4366.049999999999
$ws.Range("M122").Value = -1916.049999999999
$ws.Range("H132").Value = 1528.0869
$ws.Range("I132").Value = 797.1579
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 2391.4737
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = 138.5263
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 503
$ws.Range("I22").Value = 456
$ws.Range("K22").Value = 456
$ws.Range("M22").Value = -283
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 43511.523
$ws.Range("J138").Value = 43511.523
$ws.Range("L138").Value = 43511.523
$ws.Range("N138").Value = -53791.523
$ws.Range("H139").Value = 52567.918
$ws.Range("J139").Value = 52567.918
$ws.Range("L139").Value = 52567.918
$ws.Range("N139").Value = -62847.918
$ws.Range("H141").Value = 36340.6
$ws.Range("J141").Value = 37925.75
$ws.Range("L141").Value = 37925.75
$ws.Range("N141").Value = -48285.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 60684
$ws.Range("J18").Value = 60684
$ws.Range("L18").Value = 60684
$ws.Range("N18").Value = -61144
$ws.Range("H39").Value = 8025.5
$ws.Range("I39").Value = 8025.5
$ws.Range("K39").Value = 8025.5
$ws.Range("M39").Value = -7634.5
$ws.Range("H49").Value = 8025.5
$ws.Range("I49").Value = 8025.5
$ws.Range("K49").Value = 8025.5
$ws.Range("M49").Value = -7843.5
$ws.Range("H74").Value = 55604.332
$ws.Range("J74").Value = 55604.332
$ws.Range("L74").Value = 55604.332
$ws.Range("N74").Value = -57352.332
$ws.Range("H77").Value = 55604.332
$ws.Range("J77").Value = 55604.332
$ws.Range("L77").Value = 166812.996
$ws.Range("N77").Value = -175548.996
$ws.Range("H105").Value = 1333.2778
$ws.Range("I105").Value = 800.7778
$ws.Range("K105").Value = 800.7778
$ws.Range("M105").Value = 946.2222
$ws.Range("H132").Value = 2618.25
$ws.Range("J132").Value = 2975
$ws.Range("L132").Value = 8925
$ws.Range("N132").Value = -13985
$ws.Range("H134").Value = 2682.9656
$ws.Range("I134").Value = 2472.24
$ws.Range("K134").Value = 7416.719999999999
$ws.Range("M134").Value = -4881.719999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 292.7143
$ws.Range("I38").Value = 168.5
$ws.Range("J38").Value = 458.33334
$ws.Range("K38").Value = 505.5
$ws.Range("L38").Value = 1375.00002
$ws.Range("M38").Value = -158.5
$ws.Range("N38").Value = -2069.00002
$ws.Range("H64").Value = 4400
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -15540
$ws.Range("H67").Value = 4400
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -16872
$ws.Range("H80").Value = 3236
$ws.Range("J80").Value = 3277.5
$ws.Range("L80").Value = 9832.5
$ws.Range("N80").Value = -11704.5
$ws.Range("H81").Value = 3700
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 4666.6665
$ws.Range("K81").Value = 2400
$ws.Range("L81").Value = 13999.9995
$ws.Range("M81").Value = -1277
$ws.Range("N81").Value = -16245.9995
$ws.Range("H82").Value = 10796.941
$ws.Range("J82").Value = 11123
$ws.Range("L82").Value = 33369
$ws.Range("N82").Value = -34181
$ws.Range("H83").Value = 3236
$ws.Range("J83").Value = 3277.5
$ws.Range("L83").Value = 29497.5
$ws.Range("N83").Value = -38857.5
$ws.Range("H84").Value = 3700
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 4666.6665
$ws.Range("K84").Value = 7200
$ws.Range("L84").Value = 41999.9985
$ws.Range("M84").Value = -1584
$ws.Range("N84").Value = -53231.9985
$ws.Range("H85").Value = 10796.941
$ws.Range("J85").Value = 11123
$ws.Range("L85").Value = 33369
$ws.Range("N85").Value = -36177
$ws.Range("H87").Value = 15000
$ws.Range("J87").Value = 15000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47496
$ws.Range("H90").Value = 15000
$ws.Range("J90").Value = 15000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -147480
$ws.Range("H97").Value = 29900.572
$ws.Range("J97").Value = 51751
$ws.Range("L97").Value = 155253
$ws.Range("N97").Value = -156245
$ws.Range("H132").Value = 2201.923
$ws.Range("I132").Value = 1447
$ws.Range("K132").Value = 13023
$ws.Range("M132").Value = -10493

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 228885.11
$ws.Range("I70").Value = 668663
$ws.Range("K70").Value = 668663
$ws.Range("M70").Value = -668393
$ws.Range("H73").Value = 228885.11
$ws.Range("I73").Value = 668663
$ws.Range("K73").Value = 668663
$ws.Range("M73").Value = -667727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5749.6
$ws.Range("I7").Value = 4687
$ws.Range("K7").Value = 4687
$ws.Range("M7").Value = -4575
$ws.Range("H40").Value = 4250.1
$ws.Range("I40").Value = 4247.0625
$ws.Range("J40").Value = 4253.5713
$ws.Range("K40").Value = 4247.0625
$ws.Range("L40").Value = 4253.5713
$ws.Range("M40").Value = -4111.0625
$ws.Range("N40").Value = -4525.5713
$ws.Range("H68").Value = 2836.3333
$ws.Range("I68").Value = 2824.8
$ws.Range("J68").Value = 2859.4
$ws.Range("K68").Value = 2824.8
$ws.Range("L68").Value = 2859.4
$ws.Range("M68").Value = -2075.8
$ws.Range("N68").Value = -4357.4
$ws.Range("H71").Value = 2836.3333
$ws.Range("I71").Value = 2824.8
$ws.Range("J71").Value = 2859.4
$ws.Range("K71").Value = 14124
$ws.Range("L71").Value = 14297
$ws.Range("M71").Value = -10380
$ws.Range("N71").Value = -21785
$ws.Range("H126").Value = 5749.6
$ws.Range("I126").Value = 4687
$ws.Range("K126").Value = 14061
$ws.Range("M126").Value = -11591
$ws.Range("H134").Value = 103984
$ws.Range("J134").Value = 103984
$ws.Range("L134").Value = 103984
$ws.Range("N134").Value = -114124
$ws.Range("H136").Value = 4040.8096
$ws.Range("I136").Value = 4639
$ws.Range("J136").Value = 3243.2222
$ws.Range("K136").Value = 13917
$ws.Range("L136").Value = 9729.6666
$ws.Range("M136").Value = -11367
$ws.Range("N136").Value = -14829.6666
